# Update countries & provincias Spain
# Applies the 22-Apr-2020 18:52 data refresh to the "Pais" sheet:
#   - Ecuador's numbers grew enough to overtake Corea del Sur in the ranking
#     (rows 29/30 swap country order; Ecuador gets fresh totals, Corea del
#     Sur keeps its previous totals, just one row further down).
#   - Irak's numbers grew enough to overtake Oman (rows 69/70 swap).
#   - Aruba's numbers grew enough to overtake Bermudas (rows 143/144 swap).
#   - Estados Unidos (row 4) and Singapur (row 32) get refreshed totals
#     without changing rank/order.
#   - The "datos actualizados" timestamp moves from 18:22 to 18:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 18:52"

# --- row 4: Estados Unidos (totals refresh, order unchanged) ---
$ws.Range("B4").Value = 822572
$ws.Range("C4").Value = 3828
$ws.Range("D4").Value = 83420
$ws.Range("E4").Value = 693098
$ws.Range("F4").Value = 14016
$ws.Range("G4").Value = 736
$ws.Range("H4").Value = 46054

# --- rows 29/30: Ecuador overtakes Corea del Sur ---
$ws.Range("A29").Value = "Ecuador"
$ws.Range("B29").Value = 10850
$ws.Range("C29").Value = 452
$ws.Range("D29").Value = 1262
$ws.Range("E29").Value = 9051
$ws.Range("F29").Value = 141
$ws.Range("G29").Value = 17
$ws.Range("H29").Value = 537

$ws.Range("A30").Value = "Corea del Sur"
$ws.Range("B30").Value = 10694
$ws.Range("C30").Value = 11
$ws.Range("D30").Value = 8277
$ws.Range("E30").Value = 2179
$ws.Range("F30").Value = 55
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 238

# --- row 32: Singapur (totals refresh, order unchanged) ---
$ws.Range("D32").Value = 896
$ws.Range("E32").Value = 9233
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 12

# --- rows 69/70: Irak overtakes Oman ---
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 1631
$ws.Range("C69").Value = 29
$ws.Range("D69").Value = 1146
$ws.Range("E69").Value = 402
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 83

$ws.Range("A70").Value = "Oman"
$ws.Range("B70").Value = 1614
$ws.Range("C70").Value = 106
$ws.Range("D70").Value = 238
$ws.Range("E70").Value = 1368
$ws.Range("F70").Value = 3
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 8

# --- rows 143/144: Aruba overtakes Bermudas ---
$ws.Range("A143").Value = "Aruba"
$ws.Range("B143").Value = 100
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 68
$ws.Range("E143").Value = 30
$ws.Range("F143").Value = 4
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 2

$ws.Range("A144").Value = "Bermudas"
$ws.Range("B144").Value = 98
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 39
$ws.Range("E144").Value = 54
$ws.Range("F144").Value = 10
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 5
